# Add a new "2021" results column (M) to the table, mirroring the
# formatting of the existing "2020" column (L), then write in the
# 2021 values for each indicator row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format, font, borders, alignment) from the
# "2020" column down into the new "2021" column in one shot.
$ws.Range("L2:L10").Copy()
$ws.Range("M2:M10").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# M2 stays blank (just the formatted/bordered cell under the header row).

# Year header
$ws.Range("M3").Value = 2021

# Data rows (2021 figures)
$ws.Range("M4").Value = 952
$ws.Range("M5").Value = 10437
$ws.Range("M6").Value = 2253
$ws.Range("M7").Value = 8184
$ws.Range("M8").Value = 14020
$ws.Range("M9").Value = 5139
$ws.Range("M10").Value = 8881

# Restore the saved view state: active cell / selection at P8.
$ws.Range("P8").Select()
